$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 9,16

$arr[0,0] = 2
$arr[0,1] = 0.6666666666666666
$arr[0,2] = 0.06446533333333333
$arr[0,3] = 0.193396
$arr[0,4] = 0.02693738696927793
$arr[0,5] = 0.02693738696927793
$arr[0,6] = 3
$arr[0,7] = 1
$arr[0,8] = 0.05661333333333334
$arr[0,9] = 0.16984
$arr[0,10] = 0.0204119846136133
$arr[0,11] = 0.02041198461361329
$arr[0,12] = 0.003649597404444445
$arr[0,13] = 0.03284637664000001
$arr[0,14] = 0.0005498455283478484
$arr[0,15] = 0.0005498455283478483

$arr[1,0] = 2
$arr[1,1] = 0.6666666666666666
$arr[1,2] = 0.06446533333333333
$arr[1,3] = 0.193396
$arr[1,4] = 0.02693738696927793
$arr[1,5] = 0.02693738696927793
$arr[1,6] = 2
$arr[1,7] = 0.6666666666666666
$arr[1,8] = 0.2780386666666667
$arr[1,9] = 0.8341160000000001
$arr[1,10] = 0.100247073468963
$arr[1,11] = 0.1002470734689629
$arr[1,12] = 0.01792385532622223
$arr[1,13] = 0.161314697936
$arr[1,14] = 0.00270039421057109
$arr[1,15] = 0.00270039421057109

$arr[2,0] = 2
$arr[2,1] = 0.6666666666666666
$arr[2,2] = 0.06446533333333333
$arr[2,3] = 0.193396
$arr[2,4] = 0.02693738696927793
$arr[2,5] = 0.02693738696927793
$arr[2,6] = 3
$arr[2,7] = 1
$arr[2,8] = 2.438882
$arr[2,9] = 7.316646
$arr[2,10] = 0.8793409419174237
$arr[2,11] = 0.8793409419174237
$arr[2,12] = 0.1572233410906667
$arr[2,13] = 1.415010069816
$arr[2,14] = 0.02368714723035899
$arr[2,15] = 0.02368714723035899

$arr[3,0] = 3
$arr[3,1] = 1
$arr[3,2] = 1.843761666666667
$arr[3,3] = 5.531285
$arr[3,4] = 0.7704314695358874
$arr[3,5] = 0.7704314695358874
$arr[3,6] = 3
$arr[3,7] = 1
$arr[3,8] = 0.05661333333333334
$arr[3,9] = 0.16984
$arr[3,10] = 0.0204119846136133
$arr[3,11] = 0.02041198461361329
$arr[3,12] = 0.1043814938222222
$arr[3,13] = 0.9394334444000002
$arr[3,14] = 0.01572603530201001
$arr[3,15] = 0.01572603530201001

$arr[4,0] = 3
$arr[4,1] = 1
$arr[4,2] = 1.843761666666667
$arr[4,3] = 5.531285
$arr[4,4] = 0.7704314695358874
$arr[4,5] = 0.7704314695358874
$arr[4,6] = 2
$arr[4,7] = 0.6666666666666666
$arr[4,8] = 0.2780386666666667
$arr[4,9] = 0.8341160000000001
$arr[4,10] = 0.100247073468963
$arr[4,11] = 0.1002470734689629
$arr[4,12] = 0.5126370354511112
$arr[4,13] = 4.613733319060001
$arr[4,14] = 0.0772335001293652
$arr[4,15] = 0.07723350012936518

$arr[5,0] = 3
$arr[5,1] = 1
$arr[5,2] = 1.843761666666667
$arr[5,3] = 5.531285
$arr[5,4] = 0.7704314695358874
$arr[5,5] = 0.7704314695358874
$arr[5,6] = 3
$arr[5,7] = 1
$arr[5,8] = 2.438882
$arr[5,9] = 7.316646
$arr[5,10] = 0.8793409419174237
$arr[5,11] = 0.8793409419174237
$arr[5,12] = 4.496717141123334
$arr[5,13] = 40.47045427011
$arr[5,14] = 0.6774719341045121
$arr[5,15] = 0.6774719341045121

$arr[6,0] = 3
$arr[6,1] = 1
$arr[6,2] = 0.4849276666666666
$arr[6,3] = 1.454783
$arr[6,4] = 0.2026311434948347
$arr[6,5] = 0.2026311434948347
$arr[6,6] = 3
$arr[6,7] = 1
$arr[6,8] = 0.05661333333333334
$arr[6,9] = 0.16984
$arr[6,10] = 0.0204119846136133
$arr[6,11] = 0.02041198461361329
$arr[6,12] = 0.02745337163555556
$arr[6,13] = 0.24708034472
$arr[6,14] = 0.004136103783255434
$arr[6,15] = 0.004136103783255433

$arr[7,0] = 3
$arr[7,1] = 1
$arr[7,2] = 0.4849276666666666
$arr[7,3] = 1.454783
$arr[7,4] = 0.2026311434948347
$arr[7,5] = 0.2026311434948347
$arr[7,6] = 2
$arr[7,7] = 0.6666666666666666
$arr[7,8] = 0.2780386666666667
$arr[7,9] = 0.8341160000000001
$arr[7,10] = 0.100247073468963
$arr[7,11] = 0.1002470734689629
$arr[7,12] = 0.1348286418697778
$arr[7,13] = 1.213457776828
$arr[7,14] = 0.02031317912902667
$arr[7,15] = 0.02031317912902667

$arr[8,0] = 3
$arr[8,1] = 1
$arr[8,2] = 0.4849276666666666
$arr[8,3] = 1.454783
$arr[8,4] = 0.2026311434948347
$arr[8,5] = 0.2026311434948347
$arr[8,6] = 3
$arr[8,7] = 1
$arr[8,8] = 2.438882
$arr[8,9] = 7.316646
$arr[8,10] = 0.8793409419174237
$arr[8,11] = 0.8793409419174237
$arr[8,12] = 1.182681357535333
$arr[8,13] = 10.644132217818
$arr[8,14] = 0.1781818605825526
$arr[8,15] = 0.1781818605825526

$ws.Range("E2:T10").Value = $arr
Write-Host "done"
